# New crime data collected — weekly CompStat refresh for 24th Precinct
# (week of 3/17/2025-3/23/2025 -> 3/24/2025-3/30/2025), plus the updated
# volume/number and the full block of crime-count/percentage figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: force a cell that currently looks like (or should become) a
# number into the literal text "***.*"/"0" placeholder used throughout
# this sheet for N/A rows, while keeping the same visual style as its
# neighbours (style 13, the "dash"/placeholder style already used on
# this sheet). We do this by temporarily forcing a text number format
# (so Excel doesn't re-parse the literal "0" back into a number), then
# re-pasting the formatting (only) from a cell that already has the
# target placeholder style, restoring the normal "General" look.
# ---------------------------------------------------------------------
function Set-PlaceholderText {
    param($cellRef, $text, $formatSourceRef)
    $target = $ws.Range($cellRef)
    $target.NumberFormat = "@"
    $target.Value = $text
    $ws.Range($formatSourceRef).Copy()
    $target.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Helper: the inverse — a cell currently holding the literal placeholder
# text needs to become a real number, again re-using an existing
# numeric cell's style so the look matches its row/column neighbours.
# ---------------------------------------------------------------------
function Set-PlaceholderNumber {
    param($cellRef, $number, $formatSourceRef)
    $target = $ws.Range($cellRef)
    $target.Value = $number
    $ws.Range($formatSourceRef).Copy()
    $target.PasteSpecial(-4122)
    $target.Value = $number
}

# ---------------------------------------------------------------------
# Header text updates
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/24/2025  Through  3/30/2025"

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 7
$ws.Range("K15").Value = 133.333333333333
$ws.Range("L15").Value = 75
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = -12.5

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -10
$ws.Range("I16").Value = 28
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = -30
$ws.Range("L16").Value = -3.448275862068
$ws.Range("M16").Value = -47.169811320754
$ws.Range("N16").Value = -84.864864864864

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -88.888888888888
$ws.Range("J17").Value = 32
$ws.Range("K17").Value = 15.625
$ws.Range("L17").Value = 8.823529411764
$ws.Range("M17").Value = 32.142857142857
$ws.Range("N17").Value = -58.426966292134

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 4
Set-PlaceholderText "D18" "0" "C14"
Set-PlaceholderText "E18" "***.*" "E14"
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 275
$ws.Range("I18").Value = 39
$ws.Range("K18").Value = 30
$ws.Range("L18").Value = 8.333333333333
$ws.Range("M18").Value = 21.875
$ws.Range("N18").Value = -86.363636363636

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 83.333333333333
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 66.666666666666
$ws.Range("I19").Value = 108
$ws.Range("J19").Value = 99
$ws.Range("K19").Value = 9.090909090909
$ws.Range("L19").Value = -18.796992481203
$ws.Range("M19").Value = -8.474576271186
$ws.Range("N19").Value = -55

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
Set-PlaceholderNumber "D20" 1 "C15"
Set-PlaceholderNumber "E20" 100 "L14"
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 5
$ws.Range("I20").Value = 9
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = -30.769230769230
$ws.Range("L20").Value = -67.857142857142
$ws.Range("M20").Value = 12.5
$ws.Range("N20").Value = -95.714285714285

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = 81.818181818181
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = 40
$ws.Range("I21").Value = 228
$ws.Range("J21").Value = 217
$ws.Range("K21").Value = 5.069124423963
$ws.Range("L21").Value = -14.285714285714
$ws.Range("M21").Value = -5.785123966942
$ws.Range("N21").Value = -77.734375

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
Set-PlaceholderText "C22" "0" "C14"
$ws.Range("F22").Value = 4
Set-PlaceholderText "G22" "0" "C14"
Set-PlaceholderText "H22" "***.*" "E14"

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
Set-PlaceholderNumber "D23" 1 "C15"
Set-PlaceholderNumber "E23" 100 "L14"
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 500
$ws.Range("I23").Value = 36
$ws.Range("J23").Value = 17
$ws.Range("K23").Value = 111.764705882353
$ws.Range("L23").Value = 38.461538461538
$ws.Range("M23").Value = 38.461538461538

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 150
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = 61.290322580645
$ws.Range("I24").Value = 447
$ws.Range("J24").Value = 319
$ws.Range("K24").Value = 40.125391849529
$ws.Range("L24").Value = -4.077253218884
$ws.Range("M24").Value = 91.845493562231

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 24
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 113
$ws.Range("G25").Value = 63
$ws.Range("H25").Value = 79.365079365079
$ws.Range("I25").Value = 313
$ws.Range("J25").Value = 213
$ws.Range("K25").Value = 46.948356807511
$ws.Range("L25").Value = -2.1875

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 19
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = -24
$ws.Range("I26").Value = 72
$ws.Range("J26").Value = 70
$ws.Range("K26").Value = 2.857142857142
$ws.Range("L26").Value = 14.285714285714
$ws.Range("M26").Value = -4

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 7
$ws.Range("K27").Value = 75
$ws.Range("L27").Value = 16.666666666666

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
Set-PlaceholderNumber "C28" 3 "C15"
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 16.666666666666
$ws.Range("I28").Value = 16
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = 45.454545454545
$ws.Range("L28").Value = 6.666666666666

# ---------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------
$ws.Range("N29").Value = -75

# ---------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------
$ws.Range("N30").Value = -87.5

# ---------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------
Set-PlaceholderText "D31" "0" "C14"
Set-PlaceholderText "E31" "***.*" "E14"
